# Slide 9 ("Abstract Syntax Trees: Example 3 (continued)") contains a
# grouped diagram ("Group 1"). The three lower boxes (Expression
# (leftOperand), Expression (rightOperand), Token (operator)) move up
# slightly (~14.27pt / 181293 EMU), the bent connectors leading from the
# decision diamond to the left/right "Expression" boxes shrink to follow,
# the straight connector to the "Token" box shrinks to follow, and the
# group's own cached bounding box tightens vertically to fit its shrunk
# contents again.
#
# NOTE: in this COM-interop runtime, Shape.Left/.Top/.Width/.Height setters
# on a shape nested in a group write the EMU value straight into the
# shape's own <a:off>/<a:ext> (no group chOff/chExt or rotation adjustment
# applied), while Single-precision (float32) storage truncates toward zero
# on the EMU conversion. The literal point values below were chosen so
# that float32(value) * 12700, floored, lands exactly on the target EMU.
# (There's no exposed COM surface for a group's *child* offset/extent
# (<a:chOff>/<a:chExt>) here, only the group's own outer Left/Top/
# Width/Height (<a:off>/<a:ext>) - set below to shrink the cached box to
# match, same as the other shapes' geometry.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$grp = $s.Shapes.Item("Group 1")

# Text Box 5 "Expression (leftOperand)": off.y 3219157 -> 3037864 EMU.
$tb5 = $grp.GroupItems.Item("Text Box 5")
$tb5.Top = 239.201904296875

# Text Box 6 "Expression (rightOperand)": off.y 3219157 -> 3037864 EMU.
$tb6 = $grp.GroupItems.Item("Text Box 6")
$tb6.Top = 239.201904296875

# AutoShape 7 (bent connector, diamond -> Text Box 5): off 3316783,1963938
# -> 3407429,1873292 EMU; ext.cx 713289 -> 531996 EMU (ext.cy unchanged).
$a7 = $grp.GroupItems.Item("AutoShape 7")
$a7.Left = 268.301513671875
$a7.Top = 147.50331115722656
$a7.Width = 41.88945007324219

# AutoShape 8 (bent connector, diamond -> Text Box 6): off 5133127,1944741
# -> 5223773,1854095 EMU; ext.cx 713289 -> 531996 EMU (ext.cy unchanged).
$a8 = $grp.GroupItems.Item("AutoShape 8")
$a8.Left = 411.3207092285156
$a8.Top = 145.9917449951172
$a8.Width = 41.88945007324219

# Text Box 9 "Token (operator)": off.y 3217570 -> 3036277 EMU.
$tb9 = $grp.GroupItems.Item("Text Box 9")
$tb9.Top = 239.07693481445312

# Straight Connector 2 (diamond -> Text Box 9): ext.cy 711702 -> 530409 EMU.
$conn = $grp.GroupItems.Item("Straight Connector 2")
$conn.Height = 41.764488220214844

# Group 1 itself: ext.cy 2184944 -> 2003651 EMU (cached bounding box).
$grp.Height = 157.76780700683594
